$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.882.39"
$ws.Range("E2").Value = "  +2.58%  "
$ws.Range("D3").Value = "2.605.84"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'574.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.27%  "
$ws.Range("D6").Value = "'143.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.49%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("E8").Value = "  +0.57%  "
$ws.Range("D9").Value = "2.631.42"
$ws.Range("E9").Value = "  +1.65%  "
$ws.Range("D10").Value = "'6.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.33%  "
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("D12").Value = "'0.153"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.53%  "
$ws.Range("D13").Value = "'0.364"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.71%  "
$ws.Range("D14").Value = "3.076.42"
$ws.Range("E14").Value = "  +1.18%  "
$ws.Range("D15").Value = "60.809.58"
$ws.Range("E15").Value = "  +2.49%  "
$ws.Range("D16").Value = "'23.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("E17").Value = "  +3.14%  "
$ws.Range("D18").Value = "2.621.92"
$ws.Range("E18").Value = "  +1.26%  "
$ws.Range("D19").Value = "'11.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +9.26%  "
$ws.Range("E20").Value = "  +1.77%  "
$ws.Range("D21").Value = "'347.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.83%  "
$ws.Range("D22").Value = "'6.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.52%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("E24").Value = "  +11.90%  "
$ws.Range("D25").Value = "'63.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").Value = "'0.161"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("E28").Value = "  +3.64%  "
$ws.Range("D29").Value = "0.0₃0789"
$ws.Range("E29").Value = "  +1.54%  "
$ws.Range("E30").Value = "  +10.87%  "
$ws.Range("D31").Value = "'6.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("D33").Value = "'161.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.17%  "
$ws.Range("D34").Value = "'19.54"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.28%  "
$ws.Range("E35").Value = "  +4.35%  "
$ws.Range("D36").Value = "'0.975"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.00%  "
$ws.Range("E37").Value = "  +5.19%  "
$ws.Range("E38").Value = "  +7.86%  "
$ws.Range("D39").Value = "'37.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.42%  "
$ws.Range("D40").Value = "'3.85"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.64%  "
$ws.Range("E41").Value = "  -2.06%  "
$ws.Range("D42").Value = "'297.13"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.25%  "
$ws.Range("D43").Value = "'137.96"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.30%  "
$ws.Range("D44").Value = "'0.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("E45").Value = "  +0.95%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'19.87"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.09%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.606"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.44%  "
$ws.Range("D48").Value = "'4.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.39%  "
$ws.Range("D49").Value = "'0.0545"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.22%  "
$ws.Range("D50").Value = "'0.0241"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.56%  "
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").Value = "'10.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.87%  "
